$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows of raw session data (sessions 4, 5 and 6) ---------------
$ws.Range("B5").Value = 44894
$ws.Range("B6").Value = 44895
$ws.Range("B7").Value = 44896
$ws.Range("B5:B7").NumberFormat = "m/d/yy"

$ws.Range("C5").Value = 0.35416666666666669
$ws.Range("D5").Value = 0.5
$ws.Range("C6").Value = 0.35416666666666669
$ws.Range("D6").Value = 0.5
$ws.Range("C7").Value = 0.54166666666666663
$ws.Range("D7").Value = 0.72916666666666663
$ws.Range("C5:D7").NumberFormat = "h:mm"

# --- Extend the "temps total" (E) and "temps min" (F) formulas down to
#     row 11 so every remaining session row totals/accumulates ---------
$ws.Range("E5:E11").FormulaR1C1 = "=(RC[-1]-RC[-2])"
$ws.Range("F5:F11").FormulaR1C1 = "=R[-1]C+RC[-1]"
$ws.Range("E5:F11").NumberFormat = "h:mm;@"

# --- Selection moves to G8 after the edit ------------------------------
$ws.Range("G8").Select() | Out-Null
